$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("configurations")
$ws1.Select() | Out-Null
